# The "video views" column (E) contains some erroneously negative values
# (e.g. from an int32 overflow during a prior cleaning step). This walks
# every data row and flips any negative video-views value to its absolute
# (true, non-negative) value, leaving already-positive values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 996 }

$changed = 0
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $val = $cell.Value()
    if ($val -ne $null -and $val -lt 0) {
        $cell.Value = [Math]::Abs($val)
        $changed++
    }
}

Write-Host "max_view cleanup: fixed $changed negative video-views values"
